$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column AV: "s_n_real_pred" -------------------------------------
# Header
$ws.Range("AV1").Value = "s_n_real_pred"

# Make sure the new cells use the default "Normal" style rather than
# inheriting the customFormat applied to rows 3 & 4 (s_style 1 / 2).
$ws.Range("AV2:AV18").Style = "Normal"

# Row 2 and row 3 keep their own (non-shared) formulas; rows 4-18 share one
# formula (mirrors how the existing AT/AU "shared" columns are built).
$formula2 = "=(AQ2*X2-2*X2*(1-0.01*P2-2*0.01*AF2)/(-0.08/0.4*0.01*P2-(2*0.08/0.4+3)*0.01*AF2+1+0.08/0.4)+4*2*232000*(-0.4*0.01*AF2-0.08*0.01*P2)/((1-2*0.01*AF2)*0.4*0.08))/1000"
$formula3 = "=(AQ3*X3-2*X3*(1-0.01*P3-2*0.01*AF3)/(-0.08/0.4*0.01*P3-(2*0.08/0.4+3)*0.01*AF3+1+0.08/0.4)+4*2*232000*(-0.4*0.01*AF3-0.08*0.01*P3)/((1-2*0.01*AF3)*0.4*0.08))/1000"
$formula4 = "=(AQ4*X4-2*X4*(1-0.01*P4-2*0.01*AF4)/(-0.08/0.4*0.01*P4-(2*0.08/0.4+3)*0.01*AF4+1+0.08/0.4)+4*2*232000*(-0.4*0.01*AF4-0.08*0.01*P4)/((1-2*0.01*AF4)*0.4*0.08))/1000"

$ws.Range("AV2").Formula = $formula2
$ws.Range("AV3").Formula = $formula3
$ws.Range("AV4:AV18").Formula = $formula4

# Give the new column a width similar to its neighbors.
$ws.Columns("AV").ColumnWidth = 13.33203125

# The sheet was printed/checked through the normal page setup dialog as
# part of this edit (A4 paper).
$ws.PageSetup.PaperSize = 9

# Recalculate so every formula carries a fresh cached value.
$excel.Calculate()

# Move/extend the selection onto the newly added column, like the author
# would have done after typing the new formulas in.
$ws.Range("AV2:AV3").Select() | Out-Null
